$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","G","H","L","N","O")

$data = @(
    @(15.77687337915141, 10.61617243619246, 5.977206025344545, 11.91490136310726, 28.42984472829223, 14.00887133833256, 9.643643377518131, 17.11762359483091, 21.3190781424956),
    @(15.25067356445434, 10.47294541726249, 5.857013329162502, 11.95476457647348, 28.32502835438359, 14.04500702956846, 9.617272176399544, 17.16833313252338, 21.34790318970924),
    @(14.92033395730789, 10.38347444847146, 5.783739773013035, 11.98136935783831, 28.27155669911474, 14.0700200926753, 9.602752593859886, 17.20129952503257, 21.37139563626784),
    @(14.78410336897052, 10.34665625370583, 5.754058817645069, 11.99274588884867, 28.25251970776773, 14.08092236950588, 9.597260524988224, 17.21519480384704, 21.3824221868005),
    @(14.76139131860569, 10.34052177849254, 5.749142474568822, 11.99466724628092, 28.24952528435967, 14.08277548247049, 9.596374342772803, 17.21752998693253, 21.3843407852732),
    @(14.91850295221694, 10.38297931982001, 5.783338700213939, 11.98152062057556, 28.27128879481221, 14.07016425462229, 9.602676800681778, 17.20148505301206, 21.37153846567867),
    @(15.59706433060782, 10.56712077549443, 5.935682307775238, 11.92820417621101, 28.39145548772575, 14.02074373144033, 9.634206122515343, 17.1347288580061, 21.3278124934467),
    @(16.86192114684031, 10.9149623138469, 6.236587986678552, 11.84055744925941, 28.71256823466981, 13.94630412457744, 9.709103138431825, 17.01830393404752, 21.2881745714845),
    @(17.74127070364989, 11.16098194074414, 6.456376338537087, 11.78649313337725, 28.99906947528344, 13.90538361226044, 9.771803249283179, 16.94154015981609, 21.28730873543957),
    @(18.12872917672863, 11.27053881777045, 6.555565734803416, 11.76414566531439, 29.13997632093412, 13.88977086744018, 9.801923446510752, 16.90851188336601, 21.29306678171967),
    @(18.27351983493114, 11.3116626450596, 6.592971667560324, 11.75600673274075, 29.19481431788114, 13.88429132185038, 9.813552460950023, 16.89627616986143, 21.29613184697521),
    @(18.242424314126, 11.30282240528805, 6.584923147287594, 11.75774519886636, 29.18293886019577, 13.8854521834061, 9.811038128124322, 16.89889929068063, 21.29543239554554),
    @(18.14068056488126, 11.27392950734442, 6.558646474329523, 11.76346958281527, 29.14445847035444, 13.88931138571555, 9.802875739096981, 16.90749980885369, 21.29330122038223),
    @(18.07810458704289, 11.25618381230321, 6.542529903030231, 11.76701808609143, 29.12107952952412, 13.89173163010191, 9.79790490005248, 16.91280319616551, 21.29211100318041),
    @(17.71568557317845, 11.15377262048925, 6.449874568543279, 11.78799884472334, 28.99007025921275, 13.90646443760908, 9.769866415803079, 16.94373664735708, 21.28705624241545),
    @(17.4900420969737, 11.09032530181937, 6.392800522393347, 11.80144564554546, 28.91238093185443, 13.91627218572336, 9.753070314786308, 16.96319739836986, 21.28553120330947),
    @(17.35908203403673, 11.05361143392537, 6.35990066216104, 11.80939132691572, 28.86869425189505, 13.92219587486416, 9.743560430331774, 16.97456882607195, 21.2852331404637),
    @(17.31454329808081, 11.04114359420845, 6.348750185682188, 11.81211789712946, 28.85407532440544, 13.92425002854035, 9.740366640359412, 16.97844960931198, 21.28523166995176),
    @(17.51418488440487, 11.09710238910054, 6.398883946118261, 11.79999232400996, 28.92054807382492, 13.91519888594893, 9.754842729052548, 16.9611073386198, 21.28563360930657),
    @(18.17061851359868, 11.28242608790191, 6.566369089158784, 11.76177940700525, 29.15572128254607, 13.88816609514133, 9.8052672284142, 16.90496626838062, 21.29390319403804),
    @(18.58832991439701, 11.40141975921428, 6.674909922420401, 11.73869128422604, 29.31802347222408, 13.87302084586667, 9.839519970600838, 16.86985633290281, 21.30446350044486),
    @(18.36646155505381, 11.33811273836975, 6.617076742654825, 11.75084109802849, 29.2306270757683, 13.88087308201282, 9.821122178093029, 16.88845067799257, 21.29835573748024),
    @(17.50327376611349, 11.09403920484514, 6.396133898843523, 11.80064870120999, 28.91685265958684, 13.9156832368083, 9.754040963995536, 16.96205168396144, 21.28558550891429),
    @(16.52787020533896, 10.82243488160409, 6.155228572561881, 11.86245582447367, 28.61668917291222, 13.96402844926975, 9.687471310547052, 16.89889929068063, 21.29543239554554)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range($cols[$c] + $excelRow).Value = $rowValues[$c]
    }
}
